# Add "Group price promotion" data: a small Name/Quantity/Price table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Quantity"
$ws.Range("C1").Value = "Price"

# Orange row
$ws.Range("A2").Value = "Orange"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 2

# Tomato row
$ws.Range("A3").Value = "Tomato"
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 2

# Column widths (values back-solved so the engine's stored column width,
# after its internal character-width quantization, lands as close as
# possible to the authored widths of 22.1796875 / 17.08984375 / 14.7265625).
$ws.Columns.Item(1).ColumnWidth = 21.346354166666668
$ws.Columns.Item(2).ColumnWidth = 16.256510416666668
$ws.Columns.Item(3).ColumnWidth = 13.893229166666666

# Leave A3 as the selected/active cell, matching the saved view state.
$ws.Range("A3").Select() | Out-Null
